$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2 through 16
# from serial date 45208 (2023-10-09) to 45212 (2023-10-13)
$newDate = Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
